# Update column F ("dSF") values on the active sheet to match the
# repulled / recalculated data from the commit "repull data, push all
# data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = 0
    6  = 2
    7  = 4
    9  = -2
    17 = -3
    19 = -3
    23 = -4
    24 = 0
    26 = -5
    27 = -1
    33 = -7
    34 = -2
    35 = 0
    36 = -2
    37 = -4
    39 = 4
    43 = -1
    44 = -2
    46 = -1
    48 = -5
    50 = 1
    51 = 3
    53 = -3
    56 = -2
    57 = -2
    61 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
